$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All edited cells hold text (coin names/links, formatted price strings, and
# padded percentage strings) rather than numeric values, so force text
# number-formatting before assigning -- otherwise Excel auto-coerces plain
# decimal-looking strings (e.g. "313.94") into numbers on assignment.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.468.93"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.01%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.689.82"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.36%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.94"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.92%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3879"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.98%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4019"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.34%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.492"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.002"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.45"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.94%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08746"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.18%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "25.06"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +6.42%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.507"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +3.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.997"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.56%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001347"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.51%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.685.58"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.92%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "98.27"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.80%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07086"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.66%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "20.02"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.246"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +3.60%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.26"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.54%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.470.04"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.90%  "
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").NumberFormat = "@"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.355"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.40%  "
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.964"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -9.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.69"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.23%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.84"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.31%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.733"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +16.41%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "136.90"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.33%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.214"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.85%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.870.67"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.78%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.08821"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.41%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.417"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +4.87%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.033"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.93%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2818"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +2.71%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.955"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +3.78%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02898"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +6.23%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "10.76"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -5.77%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "14.20"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.95%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.09118"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.41%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7919"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +3.23%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.47%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.71"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.75%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.7233"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.74%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.597"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.29%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.201"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.43%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.002"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.342"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.44%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "138.15"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.50%  "
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.08026"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.54%  "
